$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44432) <- data from original row 29 (date 44274)
$ws.Cells.Item(2, 4).Value = 44274
$ws.Cells.Item(2, 10).Value = 250
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 1200
$ws.Cells.Item(2, 13).Value = 1100
$ws.Cells.Item(2, 16).Value = 1100

# Row 3 (was date 44365) <- data from original row 4 (date 44571)
$ws.Cells.Item(3, 4).Value = 44571
$ws.Cells.Item(3, 10).Value = 250
$ws.Cells.Item(3, 11).Value = 900
$ws.Cells.Item(3, 12).Value = 1000
$ws.Cells.Item(3, 13).Value = 950
$ws.Cells.Item(3, 16).Value = 950

# Row 4 (was date 44571) <- data from original row 25 (date 44523)
$ws.Cells.Item(4, 4).Value = 44523
$ws.Cells.Item(4, 10).Value = 250
$ws.Cells.Item(4, 11).Value = 1400
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = 1450
$ws.Cells.Item(4, 16).Value = 1450

# Row 5 (was date 44292) <- data from original row 2 (date 44432)
$ws.Cells.Item(5, 4).Value = 44432
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 2300
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2400
$ws.Cells.Item(5, 16).Value = 2400

# Row 6 (was date 44539) <- data from original row 8 (date 44224)
$ws.Cells.Item(6, 4).Value = 44224
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 750
$ws.Cells.Item(6, 12).Value = 800
$ws.Cells.Item(6, 13).Value = 775
$ws.Cells.Item(6, 16).Value = 775

# Row 7 (was date 44498) <- data from original row 14 (date 44376)
$ws.Cells.Item(7, 4).Value = 44376
$ws.Cells.Item(7, 10).Value = 270
$ws.Cells.Item(7, 11).Value = 2400
$ws.Cells.Item(7, 12).Value = 2500
$ws.Cells.Item(7, 13).Value = 2437
$ws.Cells.Item(7, 16).Value = 2437

# Row 8 (was date 44224) <- data from original row 12 (date 44417)
$ws.Cells.Item(8, 4).Value = 44417
$ws.Cells.Item(8, 10).Value = 250
$ws.Cells.Item(8, 11).Value = 4000
$ws.Cells.Item(8, 12).Value = 4500
$ws.Cells.Item(8, 13).Value = 4250
$ws.Cells.Item(8, 16).Value = 4250

# Row 9 (was date 44349) <- data from original row 24 (date 44250)
$ws.Cells.Item(9, 4).Value = 44250
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 12).Value = 1200
$ws.Cells.Item(9, 13).Value = 1100
$ws.Cells.Item(9, 16).Value = 1100

# Row 10 (was date 44326) <- data from original row 11 (date 44435)
$ws.Cells.Item(10, 4).Value = 44435
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 2300
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2400
$ws.Cells.Item(10, 16).Value = 2400

# Row 11 (was date 44435) <- data from original row 5 (date 44292)
$ws.Cells.Item(11, 4).Value = 44292
$ws.Cells.Item(11, 10).Value = 250
$ws.Cells.Item(11, 11).Value = 1800
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 1900
$ws.Cells.Item(11, 16).Value = 1900

# Row 12 (was date 44417) <- data from original row 26 (date 44260)
$ws.Cells.Item(12, 4).Value = 44260
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 900
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = 950
$ws.Cells.Item(12, 16).Value = 950

# Row 13 (was date 44362) <- data from original row 18 (date 44532)
$ws.Cells.Item(13, 4).Value = 44532
$ws.Cells.Item(13, 10).Value = 300
$ws.Cells.Item(13, 11).Value = 1000
$ws.Cells.Item(13, 12).Value = 1200
$ws.Cells.Item(13, 13).Value = 1100
$ws.Cells.Item(13, 16).Value = 1100

# Row 14 (was date 44376) <- data from original row 19 (date 44442)
$ws.Cells.Item(14, 4).Value = 44442
$ws.Cells.Item(14, 10).Value = 240
$ws.Cells.Item(14, 11).Value = 2300
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 13).Value = 2400
$ws.Cells.Item(14, 16).Value = 2400

# Row 15 (was date 44313) <- data from original row 16 (date 44540)
$ws.Cells.Item(15, 4).Value = 44540
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 900
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 13).Value = 950
$ws.Cells.Item(15, 16).Value = 950

# Row 16 (was date 44540) <- data from original row 27 (date 44474)
$ws.Cells.Item(16, 4).Value = 44474
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 2000
$ws.Cells.Item(16, 12).Value = 2500
$ws.Cells.Item(16, 13).Value = 2250
$ws.Cells.Item(16, 16).Value = 2250

# Row 17 (was date 44280) <- data from original row 15 (date 44313)
$ws.Cells.Item(17, 4).Value = 44313
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 900
$ws.Cells.Item(17, 12).Value = 1000
$ws.Cells.Item(17, 13).Value = 950
$ws.Cells.Item(17, 16).Value = 950

# Row 18 (was date 44532) <- data from original row 28 (date 44302)
$ws.Cells.Item(18, 4).Value = 44302
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 900
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = 950
$ws.Cells.Item(18, 16).Value = 950

# Row 19 (was date 44442) <- data from original row 6 (date 44539)
$ws.Cells.Item(19, 4).Value = 44539
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 900
$ws.Cells.Item(19, 12).Value = 1000
$ws.Cells.Item(19, 13).Value = 950
$ws.Cells.Item(19, 16).Value = 950

# Row 20 (was date 44536) <- data from original row 10 (date 44326)
$ws.Cells.Item(20, 4).Value = 44326
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 2700
$ws.Cells.Item(20, 12).Value = 2800
$ws.Cells.Item(20, 13).Value = 2750
$ws.Cells.Item(20, 16).Value = 2750

# Row 22 (was date 44330) <- data from original row 23 (date 44494)
$ws.Cells.Item(22, 4).Value = 44494
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 2400
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2450
$ws.Cells.Item(22, 16).Value = 2450

# Row 23 (was date 44494) <- data from original row 20 (date 44536)
$ws.Cells.Item(23, 4).Value = 44536
$ws.Cells.Item(23, 10).Value = 250
$ws.Cells.Item(23, 11).Value = 900
$ws.Cells.Item(23, 12).Value = 1000
$ws.Cells.Item(23, 13).Value = 950
$ws.Cells.Item(23, 16).Value = 950

# Row 24 (was date 44250) <- data from original row 17 (date 44280)
$ws.Cells.Item(24, 4).Value = 44280
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 1400
$ws.Cells.Item(24, 12).Value = 1500
$ws.Cells.Item(24, 13).Value = 1450
$ws.Cells.Item(24, 16).Value = 1450

# Row 25 (was date 44523) <- data from original row 13 (date 44362)
$ws.Cells.Item(25, 4).Value = 44362
$ws.Cells.Item(25, 10).Value = 250
$ws.Cells.Item(25, 11).Value = 2800
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 13).Value = 2900
$ws.Cells.Item(25, 16).Value = 2900

# Row 26 (was date 44260) <- data from original row 3 (date 44365)
$ws.Cells.Item(26, 4).Value = 44365
$ws.Cells.Item(26, 10).Value = 250
$ws.Cells.Item(26, 11).Value = 2400
$ws.Cells.Item(26, 12).Value = 2500
$ws.Cells.Item(26, 13).Value = 2450
$ws.Cells.Item(26, 16).Value = 2450

# Row 27 (was date 44474) <- data from original row 22 (date 44330)
$ws.Cells.Item(27, 4).Value = 44330
$ws.Cells.Item(27, 10).Value = 250
$ws.Cells.Item(27, 11).Value = 2800
$ws.Cells.Item(27, 12).Value = 3000
$ws.Cells.Item(27, 13).Value = 2900
$ws.Cells.Item(27, 16).Value = 2900

# Row 28 (was date 44302) <- data from original row 9 (date 44349)
$ws.Cells.Item(28, 4).Value = 44349
$ws.Cells.Item(28, 10).Value = 250
$ws.Cells.Item(28, 11).Value = 2800
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 2900
$ws.Cells.Item(28, 16).Value = 2900

# Row 29 (was date 44274) <- data from original row 7 (date 44498)
$ws.Cells.Item(29, 4).Value = 44498
$ws.Cells.Item(29, 10).Value = 270
$ws.Cells.Item(29, 11).Value = 2000
$ws.Cells.Item(29, 12).Value = 2300
$ws.Cells.Item(29, 13).Value = 2150
$ws.Cells.Item(29, 16).Value = 2150
